$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, report week dates) ---
$ws.Range("A8").Value = "Volume 32   Number  9"
$ws.Range("C9").Value = "Report Covering the Week  2/24/2025  Through  3/2/2025"

# --- Column H width bestFit tweak (matches column E width) ---
$ws.Columns.Item(8).ColumnWidth = $ws.Columns.Item(5).ColumnWidth

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -71.428571428571
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = 15.384615384615
$ws.Range("L16").Value = 50
$ws.Range("M16").Value = -46.428571428571
$ws.Range("N16").Value = -89.208633093525
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 12.5
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = 84.615384615384
$ws.Range("L17").Value = 118.181818181818
$ws.Range("M17").Value = 140
$ws.Range("N17").Value = 71.428571428571
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 83.333333333333
$ws.Range("I18").Value = 35
$ws.Range("J18").Value = 16
$ws.Range("K18").Value = 118.75
$ws.Range("L18").Value = 20.689655172413
$ws.Range("M18").Value = 66.666666666666
$ws.Range("N18").Value = -86.538461538461
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = -27.272727272727
$ws.Range("I19").Value = 58
$ws.Range("J19").Value = 62
$ws.Range("K19").Value = -6.451612903225
$ws.Range("L19").Value = -30.120481927710
$ws.Range("M19").Value = -12.121212121212
$ws.Range("N19").Value = -62.820512820512
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 6
$ws.Range("E20").Value = -83.333333333333
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 11
$ws.Range("H20").Value = -36.363636363636
$ws.Range("I20").Value = 22
$ws.Range("J20").Value = 20
$ws.Range("K20").Value = 10
$ws.Range("L20").Value = -12
$ws.Range("M20").Value = 4.761904761904
$ws.Range("N20").Value = -96.502384737678
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 18
$ws.Range("E21").Value = -38.888888888888
$ws.Range("F21").Value = 53
$ws.Range("G21").Value = 65
$ws.Range("H21").Value = -18.461538461538
$ws.Range("I21").Value = 156
$ws.Range("J21").Value = 125
$ws.Range("K21").Value = 24.8
$ws.Range("L21").Value = -1.886792452830
$ws.Range("M21").Value = 6.849315068493
$ws.Range("N21").Value = -87
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 6
$ws.Range("K22").Value = 20
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = 20
$ws.Range("C24").Value = 43
$ws.Range("D24").Value = 51
$ws.Range("E24").Value = -15.686274509803
$ws.Range("F24").Value = 197
$ws.Range("G24").Value = 161
$ws.Range("H24").Value = 22.360248447205
$ws.Range("I24").Value = 367
$ws.Range("J24").Value = 297
$ws.Range("K24").Value = 23.569023569023
$ws.Range("L24").Value = 32.971014492753
$ws.Range("M24").Value = 198.373983739837
$ws.Range("C25").Value = 35
$ws.Range("D25").Value = 42
$ws.Range("E25").Value = -16.666666666666
$ws.Range("F25").Value = 159
$ws.Range("G25").Value = 129
$ws.Range("H25").Value = 23.255813953488
$ws.Range("I25").Value = 302
$ws.Range("J25").Value = 234
$ws.Range("K25").Value = 29.059829059829
$ws.Range("L25").Value = 54.871794871794
$ws.Range("C26").Value = 8
$ws.Range("E26").Value = -11.111111111111
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 12
$ws.Range("I26").Value = 58
$ws.Range("J26").Value = 41
$ws.Range("K26").Value = 41.463414634146
$ws.Range("L26").Value = 75.757575757575
$ws.Range("M26").Value = 31.818181818181
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 4
$ws.Range("L28").Value = 60
$ws.Range("F31").Value = 3
$ws.Range("I31").Value = 3
$ws.Range("L31").Value = 200

# --- Cells converting text-placeholder -> numeric value (style 15, #,##0) ---
$ws.Range("D28").Value = 1
$ws.Range("I15").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("C31").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("D31").Value = 1
$ws.Range("I15").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("G31").Value = 1
$ws.Range("I15").Copy()
$ws.Range("G31").PasteSpecial(-4122)
$ws.Range("J31").Value = 1
$ws.Range("I15").Copy()
$ws.Range("J31").PasteSpecial(-4122)

# --- Cells converting text-placeholder -> numeric value (style 14, % format) ---
$ws.Range("E28").Value = 100
$ws.Range("N14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E31").Value = 0
$ws.Range("N14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("H31").Value = 200
$ws.Range("N14").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("K31").Value = 200
$ws.Range("N14").Copy()
$ws.Range("K31").PasteSpecial(-4122)

# --- Cells converting numeric -> text placeholder "0" (style 13, forced text) ---
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("A14").Copy()
$ws.Range("D22").PasteSpecial(-4122)

# --- Cells converting numeric -> text placeholder "***.*" (style 13, plain text is non-numeric already) ---
$ws.Range("A14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"

$excel.CutCopyMode = $false

